$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 8
$ws.Range("B8").Value = "MMSeqs2_97"
$ws.Range("D8").Value = 0.43
$ws.Range("E8").Value = 0.87
$ws.Range("F8").Value = 0.46
$ws.Range("G8").Value = 0.6
$ws.Range("H8").Value = 0.74

# Row 9
$ws.Range("B9").Value = "MMSeqs2_97"
$ws.Range("D9").Value = 0.31
$ws.Range("E9").Value = 0.85
$ws.Range("G9").Value = 0.48
$ws.Range("H9").Value = 0.65

# Row 10
$ws.Range("B10").Value = "MMSeqs2_97"
$ws.Range("D10").Value = 0.16
$ws.Range("E10").Value = 0.87
$ws.Range("F10").Value = 0.16
$ws.Range("G10").Value = 0.28
$ws.Range("H10").Value = 0.47
